$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated Betfair back/lay odds values for 2026-01-14 workbook.
# Each assignment below mirrors one changed cell from the source diff.

# Row 2
$ws.Range("L2").Value = 1.46
$ws.Range("O2").Value = 1.38
$ws.Range("T2").Value = 1.88
$ws.Range("X2").Value = 12.5
$ws.Range("AA2").Value = 32
$ws.Range("AE2").Value = 980
$ws.Range("AF2").Value = 980
$ws.Range("AG2").Value = 970
$ws.Range("AH2").Value = 970
$ws.Range("AK2").Value = 65
$ws.Range("AO2").Value = 1000

# Row 3
$ws.Range("G3").Value = 1.58
$ws.Range("H3").Value = 7.2
$ws.Range("P3").Value = 1.9
$ws.Range("T3").Value = 1.04
$ws.Range("W3").Value = 2.72

# Row 4
$ws.Range("J4").Value = 3.7
$ws.Range("K4").Value = 3.75
$ws.Range("N4").Value = 3.95
$ws.Range("Q4").Value = 1.97
$ws.Range("V4").Value = 1.3
$ws.Range("W4").Value = 1.98

# Row 5
$ws.Range("H5").Value = 13
$ws.Range("J5").Value = 5.6
$ws.Range("N5").Value = 3.95
$ws.Range("T5").Value = 2.46
$ws.Range("W5").Value = 3.9
$ws.Range("Y5").Value = 34

# Row 6
$ws.Range("J6").Value = 3.45
$ws.Range("R6").Value = 1.38
$ws.Range("S6").Value = 3.05
$ws.Range("Z6").Value = 23
$ws.Range("AA6").Value = 48
$ws.Range("AE6").Value = 36
$ws.Range("AI6").Value = 48
$ws.Range("AL6").Value = 50
$ws.Range("AM6").Value = 100
$ws.Range("AO6").Value = 27

# Row 7
$ws.Range("F7").Value = 1.57
$ws.Range("H7").Value = 5.5
$ws.Range("I7").Value = 6.6
$ws.Range("J7").Value = 4.1
$ws.Range("N7").Value = 4.2
$ws.Range("R7").Value = 1.43
$ws.Range("U7").Value = 2
$ws.Range("W7").Value = 2.42
$ws.Range("Y7").Value = 27
$ws.Range("Z7").Value = 65
$ws.Range("AB7").Value = 11.5
$ws.Range("AD7").Value = 24
$ws.Range("AF7").Value = 12.5
$ws.Range("AG7").Value = 12.5
$ws.Range("AH7").Value = 26
$ws.Range("AJ7").Value = 19
$ws.Range("AK7").Value = 20
$ws.Range("AL7").Value = 40
$ws.Range("AN7").Value = 10

# Row 8
$ws.Range("L8").Value = 1.25
$ws.Range("T8").Value = 1.73
$ws.Range("W8").Value = 2.84
$ws.Range("Y8").Value = 36
$ws.Range("Z8").Value = 75
$ws.Range("AB8").Value = 13
$ws.Range("AC8").Value = 14
$ws.Range("AD8").Value = 34
$ws.Range("AF8").Value = 13
$ws.Range("AG8").Value = 12.5
$ws.Range("AH8").Value = 25
$ws.Range("AI8").Value = 95
$ws.Range("AJ8").Value = 17
$ws.Range("AK8").Value = 18
$ws.Range("AL8").Value = 36
$ws.Range("AN8").Value = 7.2

# Row 9
$ws.Range("G9").Value = 2.26
$ws.Range("I9").Value = 3.65
$ws.Range("L9").Value = 1.24
$ws.Range("S9").Value = 2.22
$ws.Range("U9").Value = 2.52
$ws.Range("W9").Value = 1.79
$ws.Range("X9").Value = 28
$ws.Range("Y9").Value = 23
$ws.Range("AH9").Value = 18
$ws.Range("AI9").Value = 46
$ws.Range("AO9").Value = 29

# Row 10
$ws.Range("H10").Value = 4.8
$ws.Range("I10").Value = 4.9
$ws.Range("K10").Value = 4.2
$ws.Range("N10").Value = 5.5
$ws.Range("AA10").Value = 100
$ws.Range("AB10").Value = 12
$ws.Range("AD10").Value = 18
$ws.Range("AF10").Value = 13
$ws.Range("AL10").Value = 26

# Row 11
$ws.Range("H11").Value = 4.1
$ws.Range("N11").Value = 5.4
$ws.Range("AA11").Value = 75

# Row 12
$ws.Range("F12").Value = 12.5
$ws.Range("H12").Value = 1.25
$ws.Range("I12").Value = 1.26
$ws.Range("J12").Value = 8
$ws.Range("K12").Value = 8.199999999999999
$ws.Range("N12").Value = 10
$ws.Range("O12").Value = 1.1
$ws.Range("P12").Value = 3.95
$ws.Range("Q12").Value = 1.32
$ws.Range("R12").Value = 2.22
$ws.Range("S12").Value = 1.77
$ws.Range("T12").Value = 1.69
$ws.Range("U12").Value = 2.38
$ws.Range("V12").Value = 4.8
$ws.Range("Y12").Value = 17.5
$ws.Range("Z12").Value = 11.5
$ws.Range("AC12").Value = 18.5
$ws.Range("AG12").Value = 46
$ws.Range("AH12").Value = 26
$ws.Range("AK12").Value = 150
$ws.Range("AM12").Value = 990
$ws.Range("AN12").Value = 95
$ws.Range("AO12").Value = 2.9

# Row 13
$ws.Range("N13").Value = 5.9
$ws.Range("R13").Value = 1.67
$ws.Range("S13").Value = 2.42
$ws.Range("T13").Value = 2.8
$ws.Range("X13").Value = 32
$ws.Range("Z13").Value = 380
$ws.Range("AF13").Value = 6.8
$ws.Range("AI13").Value = 530
$ws.Range("AK13").Value = 16.5
$ws.Range("AM13").Value = 560
$ws.Range("AN13").Value = 3.4

# Row 14
$ws.Range("F14").Value = 1.46
$ws.Range("G14").Value = 1.48
$ws.Range("H14").Value = 9.199999999999999
$ws.Range("I14").Value = 9.800000000000001
$ws.Range("J14").Value = 4.6
$ws.Range("O14").Value = 1.3
$ws.Range("W14").Value = 3.05
$ws.Range("AA14").Value = 400
$ws.Range("AB14").Value = 7.8
$ws.Range("AH14").Value = 29
$ws.Range("AN14").Value = 8

# Row 15
$ws.Range("G15").Value = 7.8
$ws.Range("H15").Value = 1.5
$ws.Range("I15").Value = 1.52
$ws.Range("V15").Value = 2.92
$ws.Range("AK15").Value = 95

# Row 16
$ws.Range("G16").Value = 1.43
$ws.Range("I16").Value = 9.6
$ws.Range("N16").Value = 2.4
$ws.Range("P16").Value = 2.4
$ws.Range("Q16").Value = 1.64
$ws.Range("S16").Value = 2.3
$ws.Range("V16").Value = 1.11
$ws.Range("W16").Value = 3.25

# Row 17
$ws.Range("G17").Value = 3.15
$ws.Range("I17").Value = 3.15
$ws.Range("J17").Value = 2.98
$ws.Range("L17").Value = 1.33
$ws.Range("U17").Value = 2.16
$ws.Range("V17").Value = 1.48
$ws.Range("W17").Value = 1.49
$ws.Range("X17").Value = 18
$ws.Range("Y17").Value = 14.5
$ws.Range("Z17").Value = 23
$ws.Range("AA17").Value = 55
$ws.Range("AB17").Value = 14.5
$ws.Range("AC17").Value = 9.6
$ws.Range("AD17").Value = 15.5
$ws.Range("AE17").Value = 38
$ws.Range("AF17").Value = 23
$ws.Range("AG17").Value = 15.5
$ws.Range("AH17").Value = 21
$ws.Range("AI17").Value = 50
$ws.Range("AJ17").Value = 55
$ws.Range("AK17").Value = 38
$ws.Range("AL17").Value = 50
$ws.Range("AM17").Value = 110
$ws.Range("AN17").Value = 32
$ws.Range("AO17").Value = 32

# Row 18
$ws.Range("H18").Value = 7.8
$ws.Range("I18").Value = 11
$ws.Range("P18").Value = 2.38
$ws.Range("T18").Value = 1.85
$ws.Range("V18").Value = 1.1
$ws.Range("X18").Value = 28
$ws.Range("Y18").Value = 38
$ws.Range("Z18").Value = 90
$ws.Range("AB18").Value = 12
$ws.Range("AC18").Value = 14.5
$ws.Range("AD18").Value = 38
$ws.Range("AF18").Value = 11.5
$ws.Range("AG18").Value = 12.5
$ws.Range("AH18").Value = 29
$ws.Range("AJ18").Value = 15.5
$ws.Range("AK18").Value = 18
$ws.Range("AL18").Value = 40
$ws.Range("AN18").Value = 7
